$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$2b$12$/Ktr48YHd6JHARagglOFgOD4FL9c5AkIce8NiPqfORjLemkEEoscS'
$ws.Range("B3").Value = '$2b$12$FCOwp2iPtCQ8YSTiF.L9P.d0VhWr8u7dBKKn42hpsJzNNiam8TfSO'
$ws.Range("B4").Value = '$2b$12$gx7rkRUyzrgh8PDGSdhsFOvG2ncNrlxvzwmQEu4BCqR9IflkUbiIu'
$ws.Range("B5").Value = '$2b$12$9xdq1qX.d9XK3ftfEHF0SuBy8JdIxKvr0b57jG5Tw2zA.PPL64hFm'
